# Auto-generated edit script: updates market-price derived columns (H-N)
# for specific leve rows across all 8 sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$updates = @{
    "ALC" = @(
        @{ Row=40; Cells=@{ "H" = 4215.1816; "I" = 2938.6365; "K" = 2938.6365; "M" = -2763.6365 } }
        @{ Row=48; Cells=@{ "H" = 2840; "I" = 2766.6667; "J" = 2950; "K" = 8300.000100000001; "L" = 8850; "M" = -8008.000100000001; "N" = -9434 } }
        @{ Row=56; Cells=@{ "H" = 2840; "I" = 2766.6667; "J" = 2950; "K" = 8300.000100000001; "L" = 8850; "M" = -7766.000100000001; "N" = -9918 } }
        @{ Row=64; Cells=@{ "H" = 8132.5835; "I" = 6331.2; "J" = 9419.286; "K" = 6331.2; "L" = 9419.286; "M" = -6083.2; "N" = -9915.286 } }
        @{ Row=67; Cells=@{ "H" = 8132.5835; "I" = 6331.2; "J" = 9419.286; "K" = 6331.2; "L" = 9419.286; "M" = -5473.2; "N" = -11135.286 } }
        @{ Row=125; Cells=@{ "H" = 20163.75; "I" = 37769.332; "K" = 339923.988; "M" = -337463.988 } }
        @{ Row=138; Cells=@{ "H" = 3819.1333; "I" = 2197.5; "J" = 4629.95; "K" = 6592.5; "L" = 13889.85; "M" = -1452.5; "N" = -24169.85 } }
        @{ Row=141; Cells=@{ "H" = 8001.0835; "I" = 3445.889; "K" = 10337.667; "M" = -5157.667000000001 } }
    )
    "ARM" = @(
        @{ Row=32; Cells=@{ "H" = 11767037; "I" = 11767037; "K" = 11767037; "M" = -11766750 } }
        @{ Row=74; Cells=@{ "H" = 3150; "I" = 3106.25; "K" = 3106.25; "M" = -2232.25 } }
        @{ Row=77; Cells=@{ "H" = 3150; "I" = 3106.25; "K" = 15531.25; "M" = -11163.25 } }
        @{ Row=122; Cells=@{ "H" = 2036.2084; "I" = 1632.5333; "J" = 2709; "K" = 4897.5999; "L" = 8127; "M" = -2447.5999; "N" = -13027 } }
        @{ Row=138; Cells=@{ "H" = 51999.855; "J" = 51999.855; "L" = 51999.855; "N" = -62279.855 } }
    )
    "CRP" = @(
        @{ Row=16; Cells=@{ "H" = 7983.643; "I" = 10815.833; "J" = 5859.5; "K" = 10815.833; "L" = 5859.5; "M" = -10528.833; "N" = -6433.5 } }
        @{ Row=28; Cells=@{ "H" = 19233; "J" = 20771.834; "L" = 20771.834; "N" = -21261.834 } }
        @{ Row=31; Cells=@{ "H" = 2743.4285; "I" = 2641.2; "K" = 2641.2; "M" = -2346.2 } }
        @{ Row=34; Cells=@{ "H" = 2743.4285; "I" = 2641.2; "K" = 2641.2; "M" = -2439.2 } }
        @{ Row=86; Cells=@{ "H" = 104492; "I" = 154499.5; "K" = 154499.5; "M" = -153376.5 } }
        @{ Row=89; Cells=@{ "H" = 104492; "I" = 154499.5; "K" = 772497.5; "M" = -766881.5 } }
        @{ Row=113; Cells=@{ "H" = 7983.643; "I" = 10815.833; "J" = 5859.5; "K" = 10815.833; "L" = 5859.5; "M" = -8645.833000000001; "N" = -10199.5 } }
        @{ Row=134; Cells=@{ "H" = 2301.6667; "I" = 2286.697; "K" = 6860.091; "M" = -4325.091 } }
    )
    "CUL" = @(
        @{ Row=50; Cells=@{ "H" = 796.2222; "J" = 980.8570999999999; "L" = 2942.5713; "N" = -3904.5713 } }
        @{ Row=53; Cells=@{ "H" = 796.2222; "J" = 980.8570999999999; "L" = 2942.5713; "N" = -3904.5713 } }
        @{ Row=114; Cells=@{ "H" = 421; "I" = 253.16667; "J" = 512.5454999999999; "K" = 759.50001; "L" = 1537.6365; "M" = 2494.49999; "N" = -8045.6365 } }
        @{ Row=122; Cells=@{ "H" = 1433; "I" = 1149.5; "K" = 10345.5; "M" = -7895.5 } }
        @{ Row=140; Cells=@{ "H" = 11905758; "I" = 12500996; "J" = 1000; "K" = 37502988; "L" = 3000; "M" = -37497808; "N" = -13360 } }
    )
    "GSM" = @(
        @{ Row=126; Cells=@{ "H" = 4118.6; "I" = 4048.25; "J" = 4400; "K" = 12144.75; "L" = 13200; "M" = -9674.75; "N" = -18140 } }
        @{ Row=139; Cells=@{ "H" = 0; "J" = 0; "L" = 0 }; Deletes=@("N") }
    )
    "LTW" = @(
        @{ Row=7; Cells=@{ "H" = 8046.875; "I" = 8897; "K" = 8897; "M" = -8785 } }
        @{ Row=40; Cells=@{ "H" = 27104; "I" = 27104; "J" = 0; "K" = 27104; "L" = 0; "M" = -26968 }; Deletes=@("N") }
        @{ Row=126; Cells=@{ "H" = 8046.875; "I" = 8897; "K" = 26691; "M" = -24221 } }
        @{ Row=136; Cells=@{ "H" = 3758.125; "I" = 3859.6667; "K" = 11579.0001; "M" = -9029.000100000001 } }
    )
    "WVR" = @(
        @{ Row=82; Cells=@{ "H" = 74500; "J" = 74500; "L" = 74500; "N" = -75266 } }
        @{ Row=85; Cells=@{ "H" = 74500; "J" = 74500; "L" = 74500; "N" = -77152 } }
        @{ Row=107; Cells=@{ "H" = 511.4762; "I" = 549.7778; "K" = 1649.3334; "M" = 270.6666 } }
        @{ Row=132; Cells=@{ "H" = 1642.186; "I" = 1452.85; "K" = 4358.549999999999; "M" = -1828.549999999999 } }
        @{ Row=136; Cells=@{ "H" = 1516.6; "I" = 921.1591; "K" = 2763.4773; "M" = -213.4773 } }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $row = $entry.Row
        foreach ($col in $entry.Cells.Keys) {
            $ws.Range("$col$row").Value = $entry.Cells[$col]
        }
        if ($entry.ContainsKey("Deletes")) {
            foreach ($col in $entry.Deletes) {
                $ws.Range("$col$row").ClearContents()
            }
        }
    }
}

Write-Host "Applied $($updates.Values.Count) sheet groups of updates."